# Add a new test case row (RCC101) to the "Test Cases" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$row = 26
$ws.Cells.Item($row, 1).Value = "RCC101"
$ws.Cells.Item($row, 2).Value = "OBT3"
$ws.Cells.Item($row, 3).Value = "Verify sort by most recent activity "
$ws.Cells.Item($row, 4).Value = "Y"

# Match formatting of the row above (style with vertical-top + wrap text, plus fill applied)
$ws.Range("A" + $row).Style = $ws.Range("A25").Style
$ws.Range("C" + $row).Style = $ws.Range("C25").Style
$ws.Range("D" + $row).Style = $ws.Range("D25").Style
$ws.Range("E" + $row).Style = $ws.Range("E25").Style

$ws.Range("B" + $row).Borders.LineStyle = 1
$ws.Range("B" + $row).WrapText = $true
$ws.Range("B" + $row).VerticalAlignment = -4160

$ws.Activate()
$ws.Range("A26:E26").Select()
